# Scheduled-runner market data refresh for the Leve profit sheets.
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# columns (H:N) with freshly pulled marketboard values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 6239.65
$ws.Range("I86").Value = 1265.6154
$ws.Range("J86").Value = 15477.143
$ws.Range("K86").Value = 1265.6154
$ws.Range("L86").Value = 15477.143
$ws.Range("M86").Value = -142.6153999999999
$ws.Range("N86").Value = -17723.143
$ws.Range("H89").Value = 6239.65
$ws.Range("I89").Value = 1265.6154
$ws.Range("J89").Value = 15477.143
$ws.Range("K89").Value = 6328.076999999999
$ws.Range("L89").Value = 77385.715
$ws.Range("M89").Value = -712.0769999999993
$ws.Range("N89").Value = -88617.715
$ws.Range("H103").Value = 167070
$ws.Range("J103").Value = 566.6667
$ws.Range("L103").Value = 1700.0001
$ws.Range("N103").Value = -2872.0001
$ws.Range("H111").Value = 5651.1113
$ws.Range("I111").Value = 2342.6667
$ws.Range("J111").Value = 7305.3335
$ws.Range("K111").Value = 7028.000100000001
$ws.Range("L111").Value = 21916.0005
$ws.Range("M111").Value = -3961.000100000001
$ws.Range("N111").Value = -28050.0005
$ws.Range("H129").Value = 837.7273
$ws.Range("J129").Value = 858.5714
$ws.Range("L129").Value = 2575.7142
$ws.Range("N129").Value = -12575.7142
$ws.Range("H132").Value = 3147
$ws.Range("I132").Value = 3179.6956
$ws.Range("K132").Value = 9539.086800000001
$ws.Range("M132").Value = -7009.086800000001
$ws.Range("H135").Value = 22735430
$ws.Range("I135").Value = 807.125
$ws.Range("K135").Value = 7264.125
$ws.Range("M135").Value = -4729.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2608.3635
$ws.Range("I2").Value = 2098.5715
$ws.Range("J2").Value = 3500.5
$ws.Range("K2").Value = 2098.5715
$ws.Range("L2").Value = 3500.5
$ws.Range("M2").Value = -1985.5715
$ws.Range("N2").Value = -3726.5
$ws.Range("H32").Value = 2793.7727
$ws.Range("I32").Value = 2481.8647
$ws.Range("J32").Value = 4442.4287
$ws.Range("K32").Value = 2481.8647
$ws.Range("L32").Value = 4442.4287
$ws.Range("M32").Value = -2194.8647
$ws.Range("N32").Value = -5016.4287
$ws.Range("H45").Value = 3374.476
$ws.Range("I45").Value = 2721.7896
$ws.Range("J45").Value = 3913.652
$ws.Range("K45").Value = 2721.7896
$ws.Range("L45").Value = 3913.652
$ws.Range("M45").Value = -2344.7896
$ws.Range("N45").Value = -4667.652
$ws.Range("H61").Value = 2298.2856
$ws.Range("I61").Value = 1232
$ws.Range("J61").Value = 4102.769
$ws.Range("K61").Value = 1232
$ws.Range("L61").Value = 4102.769
$ws.Range("M61").Value = -1020
$ws.Range("N61").Value = -4526.769
$ws.Range("H97").Value = 2313.5
$ws.Range("I97").Value = 1906.875
$ws.Range("K97").Value = 1906.875
$ws.Range("M97").Value = -1410.875
$ws.Range("H110").Value = 324.16666
$ws.Range("I110").Value = 324.16666
$ws.Range("K110").Value = 324.16666
$ws.Range("M110").Value = 1720.83334
$ws.Range("H116").Value = 2608.3635
$ws.Range("I116").Value = 2098.5715
$ws.Range("J116").Value = 3500.5
$ws.Range("K116").Value = 2098.5715
$ws.Range("L116").Value = 3500.5
$ws.Range("M116").Value = 195.4285
$ws.Range("N116").Value = -8088.5
$ws.Range("H122").Value = 2329.0833
$ws.Range("I122").Value = 2393.5
$ws.Range("K122").Value = 7180.5
$ws.Range("M122").Value = -4730.5
$ws.Range("H132").Value = 12418.131
$ws.Range("I132").Value = 1291.122
$ws.Range("J132").Value = 103659.6
$ws.Range("K132").Value = 3873.366
$ws.Range("L132").Value = 310978.8
$ws.Range("M132").Value = -1343.366
$ws.Range("N132").Value = -316038.8
$ws.Range("H136").Value = 2298.2856
$ws.Range("I136").Value = 1232
$ws.Range("J136").Value = 4102.769
$ws.Range("K136").Value = 3696
$ws.Range("L136").Value = 12308.307
$ws.Range("M136").Value = -1146
$ws.Range("N136").Value = -17408.307

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2608.3635
$ws.Range("I3").Value = 2098.5715
$ws.Range("J3").Value = 3500.5
$ws.Range("K3").Value = 2098.5715
$ws.Range("L3").Value = 3500.5
$ws.Range("M3").Value = -1984.5715
$ws.Range("N3").Value = -3728.5
$ws.Range("H134").Value = 2517.2563
$ws.Range("I134").Value = 2831.75
$ws.Range("J134").Value = 1079.5714
$ws.Range("K134").Value = 8495.25
$ws.Range("L134").Value = 3238.7142
$ws.Range("M134").Value = -5960.25
$ws.Range("N134").Value = -8308.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3537.2415
$ws.Range("I31").Value = 3573.4
$ws.Range("K31").Value = 3573.4
$ws.Range("M31").Value = -3278.4
$ws.Range("H34").Value = 3537.2415
$ws.Range("I34").Value = 3573.4
$ws.Range("K34").Value = 3573.4
$ws.Range("M34").Value = -3371.4
$ws.Range("H62").Value = 5811.1665
$ws.Range("J62").Value = 5973.4
$ws.Range("L62").Value = 5973.4
$ws.Range("N62").Value = -7221.4
$ws.Range("H65").Value = 5811.1665
$ws.Range("J65").Value = 5973.4
$ws.Range("L65").Value = 29867
$ws.Range("N65").Value = -36107
$ws.Range("H122").Value = 4700.3335
$ws.Range("I122").Value = 4700.3335
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 14101.0005
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -11651.0005
$ws.Range("N122").ClearContents()
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 3842.5789
$ws.Range("I132").Value = 2533.1333
$ws.Range("J132").Value = 8753
$ws.Range("K132").Value = 7599.3999
$ws.Range("L132").Value = 26259
$ws.Range("M132").Value = -5069.3999
$ws.Range("N132").Value = -31319

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 83336240
$ws.Range("I118").Value = 166668580
$ws.Range("J118").Value = 3896.6667
$ws.Range("K118").Value = 500005740
$ws.Range("L118").Value = 11690.0001
$ws.Range("M118").Value = -500004497
$ws.Range("N118").Value = -14176.0001
$ws.Range("H131").Value = 809.18475
$ws.Range("I131").Value = 804
$ws.Range("J131").Value = 809.4828
$ws.Range("K131").Value = 2412
$ws.Range("L131").Value = 2428.4484
$ws.Range("M131").Value = 2628
$ws.Range("N131").Value = -12508.4484
$ws.Range("H132").Value = 1005.3333
$ws.Range("I132").Value = 1018.625
$ws.Range("J132").Value = 899
$ws.Range("K132").Value = 9167.625
$ws.Range("L132").Value = 8091
$ws.Range("M132").Value = -6637.625
$ws.Range("N132").Value = -13151
$ws.Range("H136").Value = 1475.6428
$ws.Range("I136").Value = 1212.2307
$ws.Range("J136").Value = 4900
$ws.Range("K136").Value = 3636.6921
$ws.Range("L136").Value = 14700
$ws.Range("M136").Value = 1463.3079
$ws.Range("N136").Value = -24900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 21440.629
$ws.Range("I132").Value = 3229.625
$ws.Range("J132").Value = 47929.363
$ws.Range("K132").Value = 9688.875
$ws.Range("L132").Value = 143788.089
$ws.Range("M132").Value = -7158.875
$ws.Range("N132").Value = -148848.089

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 10001
$ws.Range("I22").Value = 10001
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 10001
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -9706
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 10001
$ws.Range("I27").Value = 10001
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 10001
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -9894
$ws.Range("N27").ClearContents()
$ws.Range("H82").Value = 3276.5
$ws.Range("I82").Value = 3253.1428
$ws.Range("J82").Value = 3331
$ws.Range("K82").Value = 3253.1428
$ws.Range("L82").Value = 3331
$ws.Range("M82").Value = -2892.1428
$ws.Range("N82").Value = -4053
$ws.Range("H85").Value = 3276.5
$ws.Range("I85").Value = 3253.1428
$ws.Range("J85").Value = 3331
$ws.Range("K85").Value = 3253.1428
$ws.Range("L85").Value = 3331
$ws.Range("M85").Value = -2005.1428
$ws.Range("N85").Value = -5827
$ws.Range("H93").Value = 3814.1
$ws.Range("I93").Value = 3804.5557
$ws.Range("K93").Value = 3804.5557
$ws.Range("M93").Value = -2556.5557
$ws.Range("H122").Value = 2454873.5
$ws.Range("I122").Value = 3270881.2
$ws.Range("J122").Value = 6850
$ws.Range("K122").Value = 9812643.600000001
$ws.Range("L122").Value = 20550
$ws.Range("M122").Value = -9810193.600000001
$ws.Range("N122").Value = -25450
$ws.Range("H132").Value = 2483.6155
$ws.Range("I132").Value = 1754.8572
$ws.Range("J132").Value = 3333.8333
$ws.Range("K132").Value = 5264.571599999999
$ws.Range("L132").Value = 10001.4999
$ws.Range("M132").Value = -2734.571599999999
$ws.Range("N132").Value = -15061.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 4133651.2
$ws.Range("I107").Value = 1033.3334
$ws.Range("K107").Value = 3100.0002
$ws.Range("M107").Value = -1180.0002
$ws.Range("H113").Value = 3379811
$ws.Range("I113").Value = 1637.1428
$ws.Range("J113").Value = 27027028
$ws.Range("K113").Value = 4911.428400000001
$ws.Range("L113").Value = 81081084
$ws.Range("M113").Value = -2741.428400000001
$ws.Range("N113").Value = -81085424
$ws.Range("H122").Value = 2152.9333
$ws.Range("I122").Value = 1950
$ws.Range("J122").Value = 3472
$ws.Range("K122").Value = 5850
$ws.Range("L122").Value = 10416
$ws.Range("M122").Value = -3400
$ws.Range("N122").Value = -15316
$ws.Range("H132").Value = 1529.1666
$ws.Range("I132").Value = 705.7778
$ws.Range("K132").Value = 2117.3334
$ws.Range("M132").Value = 412.6666
